$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Update hardcoded input values (Methalox / main table)
$ws.Range("B2").Value = 25
$ws.Range("C2").Value = 352
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 20
$ws.Range("B4").Value = 3000
$ws.Range("C5").Value = 0.001

# Update hardcoded input values (second table, rows 12-15)
$ws.Range("B12").Value = 80
$ws.Range("H12").Value = 492.8
$ws.Range("H13").Value = 110
$ws.Range("H14").Value = 22
$ws.Range("H15").Value = 0.001

# Update selected cell
$ws.Range("C4").Select()
